$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 192.2258
$ws.Range("I9").Value = 186
$ws.Range("J9").Value = 213.57143
$ws.Range("K9").Value = 186
$ws.Range("L9").Value = 213.57143
$ws.Range("M9").Value = -17
$ws.Range("N9").Value = -551.57143
$ws.Range("H62").Value = 6869
$ws.Range("J62").Value = 7493.6665
$ws.Range("L62").Value = 7493.6665
$ws.Range("N62").Value = -8741.666499999999
$ws.Range("H65").Value = 6869
$ws.Range("J65").Value = 7493.6665
$ws.Range("L65").Value = 37468.3325
$ws.Range("N65").Value = -43708.3325
$ws.Range("H92").Value = 16667781
$ws.Range("I92").Value = 19231862
$ws.Range("K92").Value = 19231862
$ws.Range("M92").Value = -19230614
$ws.Range("H106").Value = 5130471
$ws.Range("I106").Value = 6063120.5
$ws.Range("J106").Value = 900
$ws.Range("K106").Value = 6063120.5
$ws.Range("L106").Value = 900
$ws.Range("M106").Value = -6062489.5
$ws.Range("N106").Value = -2162
$ws.Range("H125").Value = 2596.2856
$ws.Range("I125").Value = 666.3333
$ws.Range("J125").Value = 3122.6365
$ws.Range("K125").Value = 5996.9997
$ws.Range("L125").Value = 28103.7285
$ws.Range("M125").Value = -3536.9997
$ws.Range("N125").Value = -33023.7285
$ws.Range("H138").Value = 2975.9395
$ws.Range("I138").Value = 1925.6923
$ws.Range("J138").Value = 3350
$ws.Range("K138").Value = 5777.0769
$ws.Range("L138").Value = 10050
$ws.Range("M138").Value = -637.0769
$ws.Range("N138").Value = -20330

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21587630
$ws.Range("I32").Value = 25109490
$ws.Range("K32").Value = 25109490
$ws.Range("M32").Value = -25109203
$ws.Range("H45").Value = 2839.0857
$ws.Range("I45").Value = 2416.75
$ws.Range("J45").Value = 4528.4287
$ws.Range("K45").Value = 2416.75
$ws.Range("L45").Value = 4528.4287
$ws.Range("M45").Value = -2039.75
$ws.Range("N45").Value = -5282.4287
$ws.Range("H61").Value = 3356.3914
$ws.Range("I61").Value = 3019.8
$ws.Range("J61").Value = 3615.3076
$ws.Range("K61").Value = 3019.8
$ws.Range("L61").Value = 3615.3076
$ws.Range("M61").Value = -2807.8
$ws.Range("N61").Value = -4039.3076
$ws.Range("H74").Value = 2720.2307
$ws.Range("I74").Value = 2397.0588
$ws.Range("J74").Value = 4917.8
$ws.Range("K74").Value = 2397.0588
$ws.Range("L74").Value = 4917.8
$ws.Range("M74").Value = -1523.0588
$ws.Range("N74").Value = -6665.8
$ws.Range("H77").Value = 2720.2307
$ws.Range("I77").Value = 2397.0588
$ws.Range("J77").Value = 4917.8
$ws.Range("K77").Value = 11985.294
$ws.Range("L77").Value = 24589
$ws.Range("M77").Value = -7617.293999999998
$ws.Range("N77").Value = -33325
$ws.Range("H109").Value = 67839.336
$ws.Range("J109").Value = 67839.336
$ws.Range("L109").Value = 67839.336
$ws.Range("N109").Value = -70613.336
$ws.Range("H110").Value = 3500
$ws.Range("I110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("M110").Value = ""
$ws.Range("H132").Value = 4551.1577
$ws.Range("I132").Value = 4696
$ws.Range("K132").Value = 14088
$ws.Range("M132").Value = -11558
$ws.Range("H136").Value = 3356.3914
$ws.Range("I136").Value = 3019.8
$ws.Range("J136").Value = 3615.3076
$ws.Range("K136").Value = 9059.400000000001
$ws.Range("L136").Value = 10845.9228
$ws.Range("M136").Value = -6509.400000000001
$ws.Range("N136").Value = -15945.9228
$ws.Range("H139").Value = 78349.7
$ws.Range("J139").Value = 80187.125
$ws.Range("L139").Value = 80187.125
$ws.Range("N139").Value = -90467.125

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1116.3667
$ws.Range("I94").Value = 476.55
$ws.Range("K94").Value = 476.55
$ws.Range("M94").Value = -25.55000000000001
$ws.Range("H134").Value = 5499825.5
$ws.Range("I134").Value = 7939914.5
$ws.Range("J134").Value = 9625
$ws.Range("K134").Value = 23819743.5
$ws.Range("L134").Value = 28875
$ws.Range("M134").Value = -23817208.5
$ws.Range("N134").Value = -33945

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1100
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 1100
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 1100
$ws.Range("M16").Value = ""
$ws.Range("N16").Value = -1674
$ws.Range("H18").Value = 25994
$ws.Range("J18").Value = 25994
$ws.Range("L18").Value = 25994
$ws.Range("N18").Value = -26454
$ws.Range("H31").Value = 8196.84
$ws.Range("I31").Value = 1802.6666
$ws.Range("K31").Value = 1802.6666
$ws.Range("M31").Value = -1507.6666
$ws.Range("H34").Value = 8196.84
$ws.Range("I34").Value = 1802.6666
$ws.Range("K34").Value = 1802.6666
$ws.Range("M34").Value = -1600.6666
$ws.Range("H98").Value = 45025.168
$ws.Range("J98").Value = 52030.2
$ws.Range("L98").Value = 52030.2
$ws.Range("N98").Value = -56522.2
$ws.Range("H105").Value = 1227.5
$ws.Range("I105").Value = 1052.6923
$ws.Range("J105").Value = 3500
$ws.Range("K105").Value = 1052.6923
$ws.Range("L105").Value = 3500
$ws.Range("M105").Value = 694.3077000000001
$ws.Range("N105").Value = -6994
$ws.Range("H113").Value = 1100
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1100
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 1100
$ws.Range("M113").Value = ""
$ws.Range("N113").Value = -5440
$ws.Range("H122").Value = 5005942.5
$ws.Range("I122").Value = 7148207
$ws.Range("K122").Value = 21444621
$ws.Range("M122").Value = -21442171
$ws.Range("H132").Value = 3242.6
$ws.Range("I132").Value = 2919.9333
$ws.Range("J132").Value = 5178.6
$ws.Range("K132").Value = 8759.7999
$ws.Range("L132").Value = 15535.8
$ws.Range("M132").Value = -6229.7999
$ws.Range("N132").Value = -20595.8
$ws.Range("H134").Value = 2070.3684
$ws.Range("I134").Value = 1774.5312
$ws.Range("K134").Value = 5323.5936
$ws.Range("M134").Value = -2788.5936

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 5797.6
$ws.Range("I18").Value = 4989
$ws.Range("J18").Value = 5999.75
$ws.Range("K18").Value = 14967
$ws.Range("L18").Value = 17999.25
$ws.Range("M18").Value = -14798
$ws.Range("N18").Value = -18337.25
$ws.Range("H34").Value = 431.66666
$ws.Range("H46").Value = 2552.125
$ws.Range("I46").Value = 424
$ws.Range("J46").Value = 3261.5
$ws.Range("K46").Value = 1272
$ws.Range("L46").Value = 9784.5
$ws.Range("M46").Value = -1181
$ws.Range("N46").Value = -9966.5
$ws.Range("H55").Value = 1221
$ws.Range("I55").Value = 1221
$ws.Range("K55").Value = 3663
$ws.Range("M55").Value = -3486
$ws.Range("H108").Value = 799.1667
$ws.Range("I108").Value = 359
$ws.Range("K108").Value = 1077
$ws.Range("M108").Value = 1803
$ws.Range("H132").Value = 1110.3889
$ws.Range("I132").Value = 623.375
$ws.Range("K132").Value = 5610.375
$ws.Range("M132").Value = -3080.375

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 9374.125
$ws.Range("I132").Value = 11748.5
$ws.Range("J132").Value = 6999.75
$ws.Range("K132").Value = 35245.5
$ws.Range("L132").Value = 20999.25
$ws.Range("M132").Value = -32715.5
$ws.Range("N132").Value = -26059.25

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 130000
$ws.Range("J36").Value = 130000
$ws.Range("L36").Value = 130000
$ws.Range("N36").Value = -131124
$ws.Range("H59").Value = 68535.875
$ws.Range("J59").Value = 68535.875
$ws.Range("L59").Value = 68535.875
$ws.Range("N59").Value = -69843.875
$ws.Range("H124").Value = 29000
$ws.Range("J124").Value = 29000
$ws.Range("L124").Value = 29000
$ws.Range("N124").Value = -38820
$ws.Range("H132").Value = 6577
$ws.Range("I132").Value = 5999.25
$ws.Range("K132").Value = 17997.75
$ws.Range("M132").Value = -15467.75
$ws.Range("H136").Value = 3299.75
$ws.Range("I136").Value = 2099.75
$ws.Range("J136").Value = 4499.75
$ws.Range("K136").Value = 6299.25
$ws.Range("L136").Value = 13499.25
$ws.Range("M136").Value = -3749.25
$ws.Range("N136").Value = -18599.25
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").Value = ""

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1003
$ws.Range("J113").Value = 1003
$ws.Range("L113").Value = 3009
$ws.Range("N113").Value = -7349
$ws.Range("H126").Value = 5811
$ws.Range("I126").Value = 6355.4287
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 19066.2861
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -16596.2861
$ws.Range("N126").Value = -10940
$ws.Range("H132").Value = 2549.7368
$ws.Range("I132").Value = 2326.818
$ws.Range("K132").Value = 6980.454000000001
$ws.Range("M132").Value = -4450.454000000001
$ws.Range("H135").Value = 30000
$ws.Range("J135").Value = 30000
$ws.Range("L135").Value = 30000
$ws.Range("N135").Value = -40140
$ws.Range("H139").Value = 77215.78
$ws.Range("J139").Value = 77215.78
$ws.Range("L139").Value = 77215.78
$ws.Range("N139").Value = -87495.78
